# Append: 2025-10-24 06:27 JST
# - refresh "取得日時" timestamps for the rows that remain
# - replace the 4 oldest listings (rows 4-8) with newly scraped listings
# - drop the old rows 9-17 (the sheet now only tracks the newest batch)
# - shrink column B and D widths slightly
# - keep hyperlinks on column F in sync with the new URLs

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newTimestamp = "2025-10-24 06:27:19"

# ------------------------------------------------------------------
# 1. Remove all existing hyperlinks up front (the engine's Hyperlinks
#    collection Delete() operates on the whole sheet, so do this before
#    touching rows/cells and re-add the ones we still need afterwards).
# ------------------------------------------------------------------
$ws.Hyperlinks.Delete()

# ------------------------------------------------------------------
# 2. Drop the trailing rows (9-17) that are no longer present.
# ------------------------------------------------------------------
$ws.Range("A9:H17").EntireRow.Delete()

# ------------------------------------------------------------------
# 3. Refresh the "取得日時" timestamp column for every remaining data row.
# ------------------------------------------------------------------
for ($row = 2; $row -le 8; $row++) {
    $ws.Cells.Item($row, 1).Value = $newTimestamp
}

# ------------------------------------------------------------------
# 4. Replace the listing details for rows 4-8 with the newly scraped data.
# ------------------------------------------------------------------

# Row 4
$ws.Cells.Item(4, 2).Value = "Javaプログラミング研修の演習サポート講師業務【経験不問】(再掲)"
$ws.Cells.Item(4, 3).Value = "システム開発"
$ws.Cells.Item(4, 4).Value = "300,000 円 ~ 500,000 円 / 固定"
$ws.Cells.Item(4, 5).Value = "期限情報なし"
$ws.Cells.Item(4, 6).Value = "https://www.lancers.jp/work/detail/5419636"
$ws.Cells.Item(4, 7).Value = 85
$ws.Cells.Item(4, 8).Value = "★Java"

# Row 5
$ws.Cells.Item(5, 2).Value = "IB報酬を得るための高性能EA開発依頼"
$ws.Cells.Item(5, 3).Value = "システム開発"
$ws.Cells.Item(5, 4).Value = "100,000 円 ~ 200,000 円 / 固定"
$ws.Cells.Item(5, 5).Value = "期限情報なし"
$ws.Cells.Item(5, 6).Value = "https://www.lancers.jp/work/detail/5419587"
$ws.Cells.Item(5, 7).Value = 68
$ws.Cells.Item(5, 8).Value = "◆開発"

# Row 6
$ws.Cells.Item(6, 2).Value = "クラウド(AWS/Azure) 運用管理 研修の演習サポート講師業務【経験不問】(再掲)"
$ws.Cells.Item(6, 3).Value = "システム開発"
$ws.Cells.Item(6, 4).Value = "200,000 円 ~ 300,000 円 / 固定"
$ws.Cells.Item(6, 5).Value = "期限情報なし"
$ws.Cells.Item(6, 6).Value = "https://www.lancers.jp/work/detail/5419638"
$ws.Cells.Item(6, 7).Value = 38
$ws.Cells.Item(6, 8).Value = "◇管理"

# Row 7
$ws.Cells.Item(7, 2).Value = "【緊急】ロリポップ Wordpress リダイレクトハッキング復旧依頼"
$ws.Cells.Item(7, 3).Value = "システム開発"
$ws.Cells.Item(7, 4).Value = "10,000 円 ~ 20,000 円 / 固定"
$ws.Cells.Item(7, 5).Value = "期限情報なし"
$ws.Cells.Item(7, 6).Value = "https://www.lancers.jp/work/detail/5419656"
$ws.Cells.Item(7, 7).Value = 25
$ws.Cells.Item(7, 8).Value = "○WordPress"

# Row 8
$ws.Cells.Item(8, 2).Value = "【急募】Google Play Consoleでのクローズテスト実施者募集!"
$ws.Cells.Item(8, 3).Value = "システム開発"
$ws.Cells.Item(8, 4).Value = "~ 5,000 円 / 固定"
$ws.Cells.Item(8, 5).Value = "期限情報なし"
$ws.Cells.Item(8, 6).Value = "https://www.lancers.jp/work/detail/5419425"
$ws.Cells.Item(8, 7).Value = 10
$ws.Cells.Item(8, 8).ClearContents()

# ------------------------------------------------------------------
# 5. Re-create the hyperlinks for column F (rows 2-8), matching the
#    (possibly updated) URL text now shown in each cell.
# ------------------------------------------------------------------
for ($row = 2; $row -le 8; $row++) {
    $cell = $ws.Cells.Item($row, 6)
    $ws.Hyperlinks.Add($cell, $cell.Value())
}

# ------------------------------------------------------------------
# 6. Tweak column widths: B 51 -> 50, D 30 -> 28.
#    The engine's ColumnWidth setter stores (input + 5/6) as the XML
#    width, so subtract 5/6 from the desired displayed width.
# ------------------------------------------------------------------
$ws.Columns.Item(2).ColumnWidth = 50 - 5/6
$ws.Columns.Item(4).ColumnWidth = 28 - 5/6
